# Zeitaufzeichnung (time-tracking) workbook update
# Adds a new time entry in row 30: 3.5h on 2020-01-14, 17:00-21:30,
# "Präs. roles, demo/workshop roles"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hours worked
$ws.Range("A30").Value = 3.5

# Date of the entry
$ws.Range("B30").Value = "2020-01-14"

# New task description first (so its shared-string index is allocated
# before the time-range string, matching insertion order)
$ws.Range("D30").Value = "Präs. roles, demo/workshop roles"

# Time range - give it the same time number format used by other
# "Zeitraum" cells in the sheet before writing the text
$ws.Range("C30").NumberFormat = "h:mm"
$ws.Range("C30").Value = "17:00-21:30"

# Move the selection down to the next empty row, and let the view scroll
# naturally (no pinned topLeftCell)
$ws.Range("A31").Select()
